# "Generate Report for Handback" -- mark the localization status sheets as
# handed back and fill in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns for the zh-cn and de-de worksheets.

$wb = $excel.ActiveWorkbook

$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3f18b0e0e40a697767e0879ce7f04eac0c03d510/e2e/a.md"

# ---------------------------------------------------------------------------
# Overview sheet: the Status columns (zh-cn / de-de) flip from "Ready for
# handoff" to "Handed back: in sync with en-US" for both data rows.
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# Widen the Status columns to fit the longer text.
$overview.Columns.Item(5).ColumnWidth = 29.144371396019366
$overview.Columns.Item(6).ColumnWidth = 29.144371396019366

# ---------------------------------------------------------------------------
# Per-locale detail sheets: zh-cn and de-de.
# Each has rows 2 (a.md) and 3 (b.md) with:
#   I = Latest Target File      (becomes a hyperlink to a.md, like col A)
#   J = Latest Handback File    (the generated .xlf handback file name)
#   K = Latest Handback DateTime
# ---------------------------------------------------------------------------
# After Hyperlinks.Add, the engine stamps a brand-new theme-coloured
# "Hyperlink" cell style onto the target range. Pull it back in line with
# the workbook's existing (non-theme) "HyperLink" cell style -- same
# underline + same RGB colour as column A already uses -- so every
# hyperlinked cell renders identically.
function Set-HyperlinkLook($rng) {
    $rng.Font.Color = 15570276
    $rng.Font.ThemeFont = 0
    $rng.Font.Underline = $true
}

function Update-LocaleSheet($sheetName, $xlfName, $handbackDate) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Widen the Status column (C) and the Latest Handback File column (J).
    $ws.Columns.Item(3).ColumnWidth = 29.144371396019366
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664

    # Re-point the existing "b.md" hyperlink so the new links we add for
    # column I land on rId3/rId5 with "b.md" bumped to rId4, matching the
    # order Excel assigns relationship ids in.
    foreach ($h in $ws.Hyperlinks) {
        if ($h.Range.Address() -eq '$A$3') {
            $h.Delete()
        }
    }

    $ws.Range("I2").Value = "a.md"
    $ws.Hyperlinks.Add($ws.Range("I2"), $urlA, "", "", "a.md")
    Set-HyperlinkLook $ws.Range("I2")
    $ws.Range("J2").Value = $xlfName
    $ws.Range("K2").Value = $handbackDate

    $ws.Hyperlinks.Add($ws.Range("A3"), $urlA.Replace("a.md", "b.md"), "", "", "b.md")
    Set-HyperlinkLook $ws.Range("A3")

    $ws.Range("I3").Value = "a.md"
    $ws.Hyperlinks.Add($ws.Range("I3"), $urlA, "", "", "a.md")
    Set-HyperlinkLook $ws.Range("I3")
    $ws.Range("J3").Value = $xlfName
    $ws.Range("K3").Value = $handbackDate
}

Update-LocaleSheet "zh-cn" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf" "2016-09-01 00:39:54"
Update-LocaleSheet "de-de" "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf" "2016-09-01 00:40:10"
